$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: "Default Constructer" section becomes plain "Constructer" ---
$ws.Range("A7").Value = "Constructer"
$ws.Range("B7").Value = "Creating Constructor "

# --- New row 8: startCollusionAnimation() ---
$ws.Range("A8").Value = "startCollusionAnimation();"
$ws.Range("B8").Value = "When collusion with something animate. If collusion is with Penguin get bigger and disappear. Else rotate. "

# --- New row 9: Destruct ---
$ws.Range("A9").Value = "Destruct"
$ws.Range("B9").Value = "object destroys itself"

# --- Row 11: External Outgoing header row gains a Return column + moves
#     "Communication with?" out to column G ---
$ws.Range("E11").Value = "Return"
$ws.Range("G11").Value = "Communication with?"

# --- Row 12: getValue() ---
$ws.Range("A12").Value = "getValue()"
$ws.Range("B12").Value = "returns the current value of the Item."
$ws.Range("C12").Value = "Value"
$ws.Range("E12").Value = "Value"

# --- New row 13: reportCollection() ---
$ws.Range("A13").Value = "reportCollection()"
$ws.Range("B13").Value = "reports the collection to the Popup-Score and gives him current Position and Value."
$ws.Range("C13").Value = "Value, Position "
$ws.Range("E13").Value = "Value Position "
$ws.Range("G13").Value = "Popup-Score"

# --- Row 15: External Incoming header row gets Parameters back in C, Return
#     moved to E, Communication with? moved to G ---
$ws.Range("C15").Value = "Parameters"
$ws.Range("E15").Value = "Return"
$ws.Range("G15").Value = "Communication with?"

# --- New row 16: setValue() ---
$ws.Range("A16").Value = "setValue()"
$ws.Range("B16").Value = "Gamemanager sets the value to the item"
$ws.Range("C16").Value = "Value "
$ws.Range("E16").Value = "true/false"
$ws.Range("G16").Value = "Gamemanager"

# --- Widen column B to fit the new, longer descriptions ---
# (target stored width is 95.42578125 "characters"; the host's ColumnWidth
# setter quantizes to its internal pixel grid, so feed it the value whose
# rounded pixel width lands on the closest obtainable width.)
$ws.Columns("B").ColumnWidth = 94.6

# --- Selection moved from A13 to A14 ---
$null = $ws.Range("A14").Select()
